$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 52639412
$ws.Range("J62").Value = 8665.666999999999
$ws.Range("L62").Value = 8665.666999999999
$ws.Range("N62").Value = -9913.666999999999

$ws.Range("H65").Value = 52639412
$ws.Range("J65").Value = 8665.666999999999
$ws.Range("L65").Value = 43328.335
$ws.Range("N65").Value = -49568.335

$ws.Range("H107").Value = 8907.154
$ws.Range("J107").Value = 6899.75
$ws.Range("L107").Value = 6899.75
$ws.Range("N107").Value = -10739.75

$ws.Range("H116").Value = 7829477.5
$ws.Range("J116").Value = 3999
$ws.Range("L116").Value = 3999
$ws.Range("N116").Value = -10883

$ws.Range("H132").Value = 2861821.8
$ws.Range("I132").Value = 4492.5312
$ws.Range("J132").Value = 33340000
$ws.Range("K132").Value = 13477.5936
$ws.Range("L132").Value = 100020000
$ws.Range("M132").Value = -10947.5936
$ws.Range("N132").Value = -100025060

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value = 731336.5600000001
$ws.Range("I137").Value = 1161501
$ws.Range("K137").Value = 3484503
$ws.Range("M137").Value = -3481953

$ws.Range("H138").Value = 166969.77
$ws.Range("J138").Value = 5309.35
$ws.Range("L138").Value = 15928.05
$ws.Range("N138").Value = -26208.05

$ws.Range("H141").Value = 4236.5806
$ws.Range("I141").Value = 3666.6897
$ws.Range("K141").Value = 11000.0691
$ws.Range("M141").Value = -5820.069100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6399.2
$ws.Range("I2").Value = 7535.778
$ws.Range("K2").Value = 7535.778
$ws.Range("M2").Value = -7422.778

$ws.Range("H32").Value = 20331.861
$ws.Range("I32").Value = 18084.934
$ws.Range("J32").Value = 31566.5
$ws.Range("K32").Value = 18084.934
$ws.Range("L32").Value = 31566.5
$ws.Range("M32").Value = -17797.934
$ws.Range("N32").Value = -32140.5

$ws.Range("H61").Value = 6382.1577
$ws.Range("I61").Value = 6513.8667
$ws.Range("J61").Value = 5888.25
$ws.Range("K61").Value = 6513.8667
$ws.Range("L61").Value = 5888.25
$ws.Range("M61").Value = -6301.8667
$ws.Range("N61").Value = -6312.25

$ws.Range("H102").Value = 6918.1177
$ws.Range("I102").Value = 6705.921
$ws.Range("K102").Value = 6705.921
$ws.Range("M102").Value = -5083.921

$ws.Range("H116").Value = 6399.2
$ws.Range("I116").Value = 7535.778
$ws.Range("K116").Value = 7535.778
$ws.Range("M116").Value = -5241.778

$ws.Range("H122").Value = 1880967.5
$ws.Range("I122").Value = 5678.091
$ws.Range("J122").Value = 6006604
$ws.Range("K122").Value = 17034.273
$ws.Range("L122").Value = 18019812
$ws.Range("M122").Value = -14584.273
$ws.Range("N122").Value = -18024712

$ws.Range("H132").Value = 2781.75
$ws.Range("I132").Value = 2289.6875
$ws.Range("K132").Value = 6869.0625
$ws.Range("M132").Value = -4339.0625

$ws.Range("H136").Value = 6382.1577
$ws.Range("I136").Value = 6513.8667
$ws.Range("J136").Value = 5888.25
$ws.Range("K136").Value = 19541.6001
$ws.Range("L136").Value = 17664.75
$ws.Range("M136").Value = -16991.6001
$ws.Range("N136").Value = -22764.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6399.2
$ws.Range("I3").Value = 7535.778
$ws.Range("K3").Value = 7535.778
$ws.Range("M3").Value = -7421.778

$ws.Range("H26").Value = 81208.86
$ws.Range("J26").Value = 109997
$ws.Range("L26").Value = 109997
$ws.Range("N26").Value = -110581

$ws.Range("H94").Value = 649.1539
$ws.Range("I94").Value = 678.2917
$ws.Range("K94").Value = 678.2917
$ws.Range("M94").Value = -227.2917

$ws.Range("H96").Value = 30600.2
$ws.Range("I96").Value = 30600.2
$ws.Range("K96").Value = 30600.2
$ws.Range("M96").Value = -27854.2

$ws.Range("H99").Value = 21551.895
$ws.Range("I99").Value = 24581.934
$ws.Range("J99").Value = 10189.25
$ws.Range("K99").Value = 24581.934
$ws.Range("L99").Value = 10189.25
$ws.Range("M99").Value = -23083.934
$ws.Range("N99").Value = -13185.25

$ws.Range("H107").Value = 1888.5555
$ws.Range("I107").Value = 1700
$ws.Range("K107").Value = 1700
$ws.Range("M107").Value = 220

$ws.Range("H134").Value = 2417.2703
$ws.Range("I134").Value = 1510.625
$ws.Range("K134").Value = 4531.875
$ws.Range("M134").Value = -1996.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4461.1113
$ws.Range("I31").Value = 2225
$ws.Range("J31").Value = 6250
$ws.Range("K31").Value = 2225
$ws.Range("L31").Value = 6250
$ws.Range("M31").Value = -1930
$ws.Range("N31").Value = -6840

$ws.Range("H34").Value = 4461.1113
$ws.Range("I34").Value = 2225
$ws.Range("J34").Value = 6250
$ws.Range("K34").Value = 2225
$ws.Range("L34").Value = 6250
$ws.Range("M34").Value = -2023
$ws.Range("N34").Value = -6654

$ws.Range("H58").Value = 2845.1
$ws.Range("J58").Value = 3488.5
$ws.Range("L58").Value = 3488.5
$ws.Range("N58").Value = -3894.5

$ws.Range("H105").Value = 7961.6665
$ws.Range("I105").Value = 9005
$ws.Range("K105").Value = 9005
$ws.Range("M105").Value = -7258

$ws.Range("H136").Value = 2845.1
$ws.Range("J136").Value = 3488.5
$ws.Range("L136").Value = 10465.5
$ws.Range("N136").Value = -15565.5

$ws.Range("H141").Value = 412521.28
$ws.Range("J141").Value = 545554.75
$ws.Range("L141").Value = 545554.75
$ws.Range("N141").Value = -555914.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 434.27274
$ws.Range("I29").Value = 355.6
$ws.Range("J29").Value = 499.83334
$ws.Range("K29").Value = 1066.8
$ws.Range("L29").Value = 1499.50002
$ws.Range("M29").Value = -789.8000000000002
$ws.Range("N29").Value = -2053.50002

$ws.Range("H31").Value = 300
$ws.Range("J31").Value = 300
$ws.Range("L31").Value = 900
$ws.Range("N31").Value = -1476

$ws.Range("H140").Value = 557068.8
$ws.Range("I140").Value = 557068.8
$ws.Range("K140").Value = 1671206.4
$ws.Range("M140").Value = -1666026.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 6492.077
$ws.Range("I97").Value = 7258.7646
$ws.Range("J97").Value = 1278.6
$ws.Range("K97").Value = 7258.7646
$ws.Range("L97").Value = 1278.6
$ws.Range("M97").Value = -6762.7646
$ws.Range("N97").Value = -2270.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 4502
$ws.Range("J19").Value = 4502
$ws.Range("L19").Value = 4502
$ws.Range("N19").Value = -4842

$ws.Range("H40").Value = 36147.74
$ws.Range("I40").Value = 45483
$ws.Range("K40").Value = 45483
$ws.Range("M40").Value = -45347

$ws.Range("H100").Value = 3228.353
$ws.Range("I100").Value = 1073.5
$ws.Range("J100").Value = 8400
$ws.Range("K100").Value = 1073.5
$ws.Range("L100").Value = 8400
$ws.Range("M100").Value = -532.5
$ws.Range("N100").Value = -9482

$ws.Range("H122").Value = 7521.857
$ws.Range("I122").Value = 8384.333000000001
$ws.Range("K122").Value = 25152.999
$ws.Range("M122").Value = -22702.999

$ws.Range("H132").Value = 359507.88
$ws.Range("I132").Value = 682348.25
$ws.Range("K132").Value = 2047044.75
$ws.Range("M132").Value = -2044514.75

$ws.Range("H136").Value = 6039.96
$ws.Range("I136").Value = 4208.8823
$ws.Range("K136").Value = 12626.6469
$ws.Range("M136").Value = -10076.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 171164.62
$ws.Range("I62").Value = 280282.06
$ws.Range("J62").Value = 18400.2
$ws.Range("K62").Value = 280282.06
$ws.Range("L62").Value = 18400.2
$ws.Range("M62").Value = -279658.06
$ws.Range("N62").Value = -19648.2

$ws.Range("H65").Value = 171164.62
$ws.Range("I65").Value = 280282.06
$ws.Range("J65").Value = 18400.2
$ws.Range("K65").Value = 1401410.3
$ws.Range("L65").Value = 92001
$ws.Range("M65").Value = -1398290.3
$ws.Range("N65").Value = -98241

$ws.Range("H105").Value = 23532.666
$ws.Range("J105").Value = 23532.666
$ws.Range("L105").Value = 23532.666
$ws.Range("N105").Value = -30520.666

$ws.Range("H107").Value = 8554.385
$ws.Range("I107").Value = 1301.5
$ws.Range("K107").Value = 3904.5
$ws.Range("M107").Value = -1984.5

$ws.Range("H122").Value = 5617.7715
$ws.Range("I122").Value = 4274.6113
$ws.Range("K122").Value = 12823.8339
$ws.Range("M122").Value = -10373.8339

$ws.Range("H132").Value = 13328.574
$ws.Range("I132").Value = 16273.639
$ws.Range("K132").Value = 48820.917
$ws.Range("M132").Value = -46290.917

$ws.Range("H136").Value = 315700.47
$ws.Range("I136").Value = 376815.4
$ws.Range("J136").Value = 2486.375
$ws.Range("K136").Value = 1130446.2
$ws.Range("L136").Value = 7459.125
$ws.Range("M136").Value = -1127896.2
$ws.Range("N136").Value = -12559.125
